# Auto-generated Excel COM-interop script to apply the diff changes
# to Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) of the workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 16773.334  # was 6102.2856
$ws.Range("I21").Value = 5149.5  # was 3786
$ws.Range("J21").Value = 40021  # was 20000
$ws.Range("K21").Value = 5149.5  # was 3786
$ws.Range("L21").Value = 40021  # was 20000
$ws.Range("M21").Value = -4681.5  # was -3318
$ws.Range("N21").Value = -40957  # was -20936
$ws.Range("H23").Value = 16773.334  # was 6102.2856
$ws.Range("I23").Value = 5149.5  # was 3786
$ws.Range("J23").Value = 40021  # was 20000
$ws.Range("K23").Value = 5149.5  # was 3786
$ws.Range("L23").Value = 40021  # was 20000
$ws.Range("M23").Value = -4915.5  # was -3552
$ws.Range("N23").Value = -40489  # was -20468
$ws.Range("H40").Value = 7000  # was 5366.6665
$ws.Range("I40").Value = 4000  # was 3050
$ws.Range("K40").Value = 4000  # was 3050
$ws.Range("M40").Value = -3825  # was -2875
$ws.Range("H74").Value = 4987.3335  # was 4989
$ws.Range("I74").Value = 4987.3335  # was 4989
$ws.Range("K74").Value = 4987.3335  # was 4989
$ws.Range("M74").Value = -4051.3335  # was -4053
$ws.Range("H77").Value = 4987.3335  # was 4989
$ws.Range("I77").Value = 4987.3335  # was 4989
$ws.Range("K77").Value = 24936.6675  # was 24945
$ws.Range("M77").Value = -20256.6675  # was -20265
$ws.Range("H100").Value = 8159.8  # was 4679.8
$ws.Range("I100").Value = 9949.5  # was 5599.75
$ws.Range("J100").Value = 6966.6665  # was 1000
$ws.Range("K100").Value = 9949.5  # was 5599.75
$ws.Range("L100").Value = 6966.6665  # was 1000
$ws.Range("M100").Value = -9408.5  # was -5058.75
$ws.Range("N100").Value = -8048.6665  # was -2082
$ws.Range("H132").Value = 2183  # was 2203.5334
$ws.Range("I132").Value = 1580.1666  # was 1650.2307
$ws.Range("K132").Value = 4740.4998  # was 4950.6921
$ws.Range("M132").Value = -2210.4998  # was -2420.6921

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1007  # was 0
$ws.Range("I26").Value = 1007  # was 0
$ws.Range("K26").Value = 1007  # was 0
$ws.Range("M26").Value = -677  # was None
$ws.Range("H132").Value = 4356  # was 1995.4783
$ws.Range("I132").Value = 4498.3335  # was 1544.95
$ws.Range("J132").Value = 4249.25  # was 4999
$ws.Range("K132").Value = 13495.0005  # was 4634.85
$ws.Range("L132").Value = 12747.75  # was 14997
$ws.Range("M132").Value = -10965.0005  # was -2104.85
$ws.Range("N132").Value = -17807.75  # was -20057

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3390.8333  # was 3731
$ws.Range("I134").Value = 2095.9375  # was 2299.7144
$ws.Range("K134").Value = 6287.8125  # was 6899.1432
$ws.Range("M134").Value = -3752.8125  # was -4364.1432

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1283.4  # was 1258.7
$ws.Range("I31").Value = 947.7143  # was 981.1667
$ws.Range("J31").Value = 2066.6667  # was 1675
$ws.Range("K31").Value = 947.7143  # was 981.1667
$ws.Range("L31").Value = 2066.6667  # was 1675
$ws.Range("M31").Value = -652.7143  # was -686.1667
$ws.Range("N31").Value = -2656.6667  # was -2265
$ws.Range("H34").Value = 1283.4  # was 1258.7
$ws.Range("I34").Value = 947.7143  # was 981.1667
$ws.Range("J34").Value = 2066.6667  # was 1675
$ws.Range("K34").Value = 947.7143  # was 981.1667
$ws.Range("L34").Value = 2066.6667  # was 1675
$ws.Range("M34").Value = -745.7143  # was -779.1667
$ws.Range("N34").Value = -2470.6667  # was -2079

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 10000  # was 0
$ws.Range("J68").Value = 10000  # was 0
$ws.Range("L68").Value = 30000  # was 0
$ws.Range("N68").Value = -31622  # was None
$ws.Range("H71").Value = 10000  # was 0
$ws.Range("J71").Value = 10000  # was 0
$ws.Range("L71").Value = 90000  # was 0
$ws.Range("N71").Value = -98112  # was None
$ws.Range("H92").Value = 699.2  # was 624.25
$ws.Range("J92").Value = 999.5  # was 1000
$ws.Range("L92").Value = 2998.5  # was 3000
$ws.Range("N92").Value = -5494.5  # was -5496
$ws.Range("H107").Value = 625  # was 636.2
$ws.Range("J107").Value = 625  # was 636.2
$ws.Range("L107").Value = 1875  # was 1908.6
$ws.Range("N107").Value = -5715  # was -5748.6
$ws.Range("H131").Value = 2590.353  # was 2527
$ws.Range("I131").Value = 2910  # was 2880
$ws.Range("J131").Value = 2492  # was 2456.4
$ws.Range("K131").Value = 8730  # was 8640
$ws.Range("L131").Value = 7476  # was 7369.200000000001
$ws.Range("M131").Value = -3690  # was -3600
$ws.Range("N131").Value = -17556  # was -17449.2
$ws.Range("H139").Value = 2500  # was 2866.5
$ws.Range("I139").Value = 2500  # was 1700
$ws.Range("J139").Value = 0  # was 4033
$ws.Range("K139").Value = 7500  # was 5100
$ws.Range("L139").Value = 0  # was 12099
$ws.Range("M139").Value = -2360  # was 40
$ws.Range("N139").ClearContents()  # was -22379

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13501.333  # was 2749
$ws.Range("I80").Value = 2498  # was 2749
$ws.Range("J80").Value = 19003  # was 0
$ws.Range("K80").Value = 2498  # was 2749
$ws.Range("L80").Value = 19003  # was 0
$ws.Range("M80").Value = -1500  # was -1751
$ws.Range("N80").Value = -20999  # was None
$ws.Range("H83").Value = 13501.333  # was 2749
$ws.Range("I83").Value = 2498  # was 2749
$ws.Range("J83").Value = 19003  # was 0
$ws.Range("K83").Value = 12490  # was 13745
$ws.Range("L83").Value = 95015  # was 0
$ws.Range("M83").Value = -7498  # was -8753
$ws.Range("N83").Value = -104999  # was None
$ws.Range("H102").Value = 2131.6667  # was 2214.5
$ws.Range("I102").Value = 2131.6667  # was 2214.5
$ws.Range("K102").Value = 2131.6667  # was 2214.5
$ws.Range("M102").Value = -509.6667000000002  # was -592.5
$ws.Range("H122").Value = 0  # was 4999
$ws.Range("I122").Value = 0  # was 4999
$ws.Range("K122").Value = 0  # was 14997
$ws.Range("M122").ClearContents()  # was -12547

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2900  # was 2685.7144
$ws.Range("I46").Value = 2900  # was 2685.7144
$ws.Range("K46").Value = 2900  # was 2685.7144
$ws.Range("M46").Value = -2712  # was -2497.7144
$ws.Range("H55").Value = 1067.875  # was 1307.3334
$ws.Range("I55").Value = 762  # was 1174.5
$ws.Range("K55").Value = 762  # was 1174.5
$ws.Range("M55").Value = -589  # was -1001.5
$ws.Range("H61").Value = 7287164.5  # was 6801420
$ws.Range("I61").Value = 10201201  # was 8501167
$ws.Range("K61").Value = 10201201  # was 8501167
$ws.Range("M61").Value = -10200999  # was -8500965
$ws.Range("H113").Value = 7287164.5  # was 6801420
$ws.Range("I113").Value = 10201201  # was 8501167
$ws.Range("K113").Value = 10201201  # was 8501167
$ws.Range("M113").Value = -10199031  # was -8498997
$ws.Range("H136").Value = 983.2  # was 772.3333
$ws.Range("I136").Value = 729.25  # was 643.1429000000001
$ws.Range("J136").Value = 1999  # was 1224.5
$ws.Range("K136").Value = 2187.75  # was 1929.4287
$ws.Range("L136").Value = 5997  # was 3673.5
$ws.Range("M136").Value = 362.25  # was 620.5712999999998
$ws.Range("N136").Value = -11097  # was -8773.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2073.6667  # was 3524.1428
$ws.Range("I81").Value = 2073.6667  # was 2611.5
$ws.Range("J81").Value = 0  # was 9000
$ws.Range("K81").Value = 4147.3334  # was 5223
$ws.Range("L81").Value = 0  # was 18000
$ws.Range("M81").Value = -3086.3334  # was -4162
$ws.Range("N81").ClearContents()  # was -20122
$ws.Range("H82").Value = 90000  # was 51500
$ws.Range("I82").Value = 0  # was 8000
$ws.Range("J82").Value = 90000  # was 95000
$ws.Range("K82").Value = 0  # was 8000
$ws.Range("L82").Value = 90000  # was 95000
$ws.Range("M82").ClearContents()  # was -7617
$ws.Range("N82").Value = -90766  # was -95766
$ws.Range("H84").Value = 2073.6667  # was 3524.1428
$ws.Range("I84").Value = 2073.6667  # was 2611.5
$ws.Range("J84").Value = 0  # was 9000
$ws.Range("K84").Value = 20736.667  # was 26115
$ws.Range("L84").Value = 0  # was 90000
$ws.Range("M84").Value = -15432.667  # was -20811
$ws.Range("N84").ClearContents()  # was -100608
$ws.Range("H85").Value = 90000  # was 51500
$ws.Range("I85").Value = 0  # was 8000
$ws.Range("J85").Value = 90000  # was 95000
$ws.Range("K85").Value = 0  # was 8000
$ws.Range("L85").Value = 90000  # was 95000
$ws.Range("M85").ClearContents()  # was -6674
$ws.Range("N85").Value = -92652  # was -97652
